$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Hello"
$ws.Range("B9").Value = "world"
$ws.Range("C9").Value = "2025-01-30 15:18:26"
